# Commit: "Update edit dan hapus hooks.js"
# The TagSave sheet (sheet3) gains three more saved tag rows (91-93),
# extending the used range from A1:A90 to A1:A93.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TagSave")

$newValues = @(
    "oq7XiX5LVt",
    "eTGkZta5kR",
    "VqRxXMuFgd"
)

$startRow = 91
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $newValues[$i]
    # Mirror the sheet's existing convention of ignoring the
    # "number stored as text" warning for this column.
    try {
        $cell.Errors.Item(9).Ignore = $true
    } catch {
    }
}

# Keep the ignored-error range in sync with the new used range (A1:A93),
# matching how the rest of the column was already configured.
try {
    $ws.Range("A1:A93").Errors.Item(9).Ignore = $true
} catch {
}
